# Add a new "log" / change-log flag column (P) to the meta_class attribute
# definition sheet, mirroring the existing class/field attribute columns
# (A..O). Every existing attribute row defaults to 0 (log disabled).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1 (English field key) and row 2 (Chinese description),
# formatted the same way as the rest of the header cells in those rows.
$ws.Range("P1").Value = "log"
$ws.Range("P2").Value = "改变时log"

# Data rows 3-37: every attribute row gets a default value of 0.
# The sheet alternates a light "s=3" style on every other data row
# (matching the existing zebra-style formatting already used in column B),
# starting with row 4.
for ($r = 3; $r -le 37; $r++) {
    $ws.Cells.Item($r, 16).Value = 0
    if (($r % 2) -eq 0) {
        $ws.Range("B4").Copy()
        $ws.Cells.Item($r, 16).PasteSpecial(-4122)
        $ws.Cells.Item($r, 16).Value = 0
    }
}

# Size the new column to fit its (Chinese) header text (best-fit width).
$ws.Columns.Item(16).ColumnWidth = 9.37

# Leave the selection on the newly added column, as in the saved file.
$ws.Range("P10").Select()
